$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 80: date text must stay literal text (not auto-converted to a date
# serial number), so force the cell to Text format before assigning, then
# reset the style back to Normal so no style index is left on the cell.
$ws.Range("A80").NumberFormat = "@"
$ws.Range("A80").Value = "2025/10/08"
$ws.Range("A80").Style = "Normal"

$ws.Range("B80").Value = "水"
$ws.Range("C80").Value = 21
$ws.Range("D80").Value = 16
